$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "EN_A0100"
$ws.Range("G17").Value = "EN_A0800"

$ws.Range("G18").Select()
